$p = $ppt.ActivePresentation

# Slide 1 title: consolidate the "First" + " " + "slide" runs into a single
# run. Re-assigning TextRange.Text to the exact same concatenated string is
# treated as a no-op by the writer's diffing, so first nudge it through a
# distinct placeholder value, then assign the real (already-correct) text so
# the underlying runs actually get rewritten/merged.
$s1 = $p.Slides.Item(1)
$t1 = $s1.Shapes.Item(1).TextFrame.TextRange
$t1.Text = "placeholder"
$t1.Text = "First slide"

# Slide 3 title: consolidate the "Third" + " " + "slide" runs into a single
# run, using the same nudge trick.
$s3 = $p.Slides.Item(3)
$t3 = $s3.Shapes.Item(1).TextFrame.TextRange
$t3.Text = "placeholder"
$t3.Text = "Third slide"
